$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 86: header row repeated (same strings as row 1, columns B-G)
$ws.Range("B86").Value = "fastText-crawl-300d-2M.vec"
$ws.Range("C86").Value = "refined-0.5(fastText-crawl-300d-2M.vec)"
$ws.Range("D86").Value = "refined-1(fastText-crawl-300d-2M.vec)"
$ws.Range("E86").Value = "refined-2(fastText-crawl-300d-2M.vec)"
$ws.Range("F86").Value = "refined-10(fastText-crawl-300d-2M.vec)"
$ws.Range("G86").Value = "refined-20(fastText-crawl-300d-2M.vec)"

# Row 87
$ws.Range("A87").Value = "turned-fastText(uni)"
$ws.Range("D87").Value = "Test accuracy: 0.85671"
$ws.Range("G87").Value = "Test accuracy: 0.86850"

# Row 88
$ws.Range("D88").Value = "Test accuracy: 0.86194"
$ws.Range("G88").Value = "Test accuracy: 0.86619"

# Row 89
$ws.Range("D89").Value = "Test accuracy: 0.86011"
$ws.Range("G89").Value = "Test accuracy: 0.85289"

# Row 90
$ws.Range("D90").Value = "Test accuracy: 0.84974"

# Row 91
$ws.Range("D91").Value = "Test accuracy: 0.84238"

# Update the view to match the target (topLeftCell B67, selection F90)
$ws.Range("F90").Select()
$excel.ActiveWindow.ScrollRow = 67
$excel.ActiveWindow.ScrollColumn = 2
